$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("G2").Value = 3.9
$ws.Range("H2").Value = 2.75
$ws.Range("I2").Value = 2.3
$ws.Range("L2").Value = 3.2
$ws.Range("O2").Value = 1.62
$ws.Range("P2").Value = 2.2
$ws.Range("Q2").Value = 2.18
$ws.Range("R2").Value = 1.69
$ws.Range("U2").Value = 5
$ws.Range("V2").Value = 1.17
$ws.Range("AO2").Value = 9
$ws.Range("AQ2").Value = 21
$ws.Range("AR2").Value = 23

# Row 3
$ws.Range("G3").Value = 1.8
$ws.Range("H3").Value = 3.3
$ws.Range("I3").Value = 5
$ws.Range("J3").Value = 2.5
$ws.Range("K3").Value = 1.95
$ws.Range("L3").Value = 5.5
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 7
$ws.Range("O3").Value = 1.5
$ws.Range("P3").Value = 2.5
$ws.Range("Q3").Value = 1.95
$ws.Range("R3").Value = 1.9
$ws.Range("S3").Value = 2.5
$ws.Range("T3").Value = 1.5
$ws.Range("U3").Value = 4.1
$ws.Range("V3").Value = 1.24
$ws.Range("Y3").Value = 1.57
$ws.Range("Z3").Value = 2.25
$ws.Range("AD3").Value = 7
$ws.Range("AF3").Value = 13
$ws.Range("AG3").Value = 19
$ws.Range("AI3").Value = 6.5
$ws.Range("AJ3").Value = 6.5
$ws.Range("AN3").Value = 10
$ws.Range("AO3").Value = 23
$ws.Range("AP3").Value = 17
$ws.Range("AQ3").Value = 51
$ws.Range("AR3").Value = 41

# Row 4
$ws.Range("U4").Value = 3.8
$ws.Range("V4").Value = 1.27

# Row 5
$ws.Range("S5").Value = 2.35
$ws.Range("T5").Value = 1.57

# Row 8
$ws.Range("G8").Value = 3.6
$ws.Range("H8").Value = 3
$ws.Range("J8").Value = 4.5
$ws.Range("K8").Value = 1.91
$ws.Range("M8").Value = 1.11
$ws.Range("N8").Value = 6.5
$ws.Range("O8").Value = 1.5
$ws.Range("P8").Value = 2.5
$ws.Range("Q8").Value = 2
$ws.Range("R8").Value = 1.85
$ws.Range("S8").Value = 2.6
$ws.Range("T8").Value = 1.48
$ws.Range("U8").Value = 4.2
$ws.Range("V8").Value = 1.22
$ws.Range("W8").Value = 5.5
$ws.Range("X8").Value = 1.14
$ws.Range("Y8").Value = 1.57
$ws.Range("Z8").Value = 2.25
$ws.Range("AA8").Value = 2.2
$ws.Range("AB8").Value = 1.62
$ws.Range("AC8").Value = 8
$ws.Range("AD8").Value = 17
$ws.Range("AE8").Value = 13
$ws.Range("AF8").Value = 41
$ws.Range("AG8").Value = 34
$ws.Range("AH8").Value = 51
$ws.Range("AI8").Value = 6.5
$ws.Range("AK8").Value = 19
$ws.Range("AL8").Value = 81
$ws.Range("AN8").Value = 6
$ws.Range("AO8").Value = 9
$ws.Range("AP8").Value = 10
$ws.Range("AS8").Value = 41

# Row 12
$ws.Range("G12").Value = 1.29
$ws.Range("H12").Value = 4.6
$ws.Range("I12").Value = 11.75
$ws.Range("J12").Value = 1.78
$ws.Range("K12").Value = 2.3
$ws.Range("L12").Value = 9
$ws.Range("N12").Value = 8
$ws.Range("O12").Value = 1.23
$ws.Range("P12").Value = 3.7
$ws.Range("S12").Value = 1.7
$ws.Range("T12").Value = 2.02
$ws.Range("W12").Value = 2.67
$ws.Range("X12").Value = 1.42
$ws.Range("Y12").Value = 1.37
$ws.Range("Z12").Value = 2.85
$ws.Range("AA12").Value = 2.12
$ws.Range("AB12").Value = 1.65
$ws.Range("AC12").Value = 6.3
$ws.Range("AD12").Value = 5.8
$ws.Range("AE12").Value = 8.5
$ws.Range("AF12").Value = 7.6
$ws.Range("AG12").Value = 11.25
$ws.Range("AH12").Value = 30
$ws.Range("AI12").Value = 8
$ws.Range("AJ12").Value = 9.5
$ws.Range("AK12").Value = 23
$ws.Range("AL12").Value = 120
$ws.Range("AN12").Value = 28
$ws.Range("AP12").Value = 35
$ws.Range("AQ12").Value = 450
$ws.Range("AR12").Value = 175
$ws.Range("AS12").Value = 110

# Row 13
$ws.Range("H13").Value = 3.4
$ws.Range("I13").Value = 4.2
$ws.Range("M13").Value = 1.06
$ws.Range("N13").Value = 10
$ws.Range("S13").Value = 1.93
$ws.Range("T13").Value = 1.93
$ws.Range("AI13").Value = 10
$ws.Range("AK13").Value = 15
$ws.Range("AN13").Value = 12

# Row 14
$ws.Range("G14").Value = 2.35
$ws.Range("J14").Value = 3.25
$ws.Range("AF14").Value = 23

# Row 15
$ws.Range("M15").Value = 1.07
$ws.Range("N15").Value = 9
$ws.Range("S15").Value = 2.1
$ws.Range("T15").Value = 1.73
$ws.Range("W15").Value = 3.75
$ws.Range("X15").Value = 1.29
$ws.Range("AJ15").Value = 6
$ws.Range("AP15").Value = 12

# Row 16
$ws.Range("Y16").Value = 1.36
$ws.Range("Z16").Value = 3
$ws.Range("AH16").Value = 67
$ws.Range("AJ16").Value = 8.5
$ws.Range("AM16").Value = 451
$ws.Range("AN16").Value = 6.5
$ws.Range("AO16").Value = 6
$ws.Range("AP16").Value = 8.5

# Row 17
$ws.Range("G17").Value = 1.73
$ws.Range("I17").Value = 5
$ws.Range("K17").Value = 2.2
$ws.Range("L17").Value = 5
$ws.Range("M17").Value = 1.06
$ws.Range("N17").Value = 10
$ws.Range("O17").Value = 1.33
$ws.Range("P17").Value = 3.4
$ws.Range("S17").Value = 2.05
$ws.Range("T17").Value = 1.8
$ws.Range("W17").Value = 3.5
$ws.Range("X17").Value = 1.3
$ws.Range("AA17").Value = 1.91
$ws.Range("AB17").Value = 1.8
$ws.Range("AC17").Value = 6.5
$ws.Range("AI17").Value = 9.5
$ws.Range("AK17").Value = 17
$ws.Range("AL17").Value = 51
$ws.Range("AM17").Value = 351
$ws.Range("AO17").Value = 23
$ws.Range("AP17").Value = 15
$ws.Range("AS17").Value = 41

# Row 18
$ws.Range("G18").Value = 1.29
$ws.Range("H18").Value = 5.5
$ws.Range("I18").Value = 9
$ws.Range("J18").Value = 1.73
$ws.Range("L18").Value = 9
$ws.Range("S18").Value = 1.73
$ws.Range("T18").Value = 2.08
$ws.Range("AA18").Value = 2.2
$ws.Range("AB18").Value = 1.62
$ws.Range("AC18").Value = 6.5
$ws.Range("AF18").Value = 7.5
$ws.Range("AH18").Value = 34
$ws.Range("AI18").Value = 12
$ws.Range("AJ18").Value = 11
$ws.Range("AQ18").Value = 126
$ws.Range("AS18").Value = 67

# Row 21
$ws.Range("M21").Value = 1.04
$ws.Range("N21").Value = 13
$ws.Range("O21").Value = 1.25
$ws.Range("P21").Value = 3.75
$ws.Range("S21").Value = 1.75
$ws.Range("T21").Value = 2.05
$ws.Range("W21").Value = 3
$ws.Range("X21").Value = 1.36

# Row 22
$ws.Range("G22").Value = 2.9
$ws.Range("H22").Value = 3.5
$ws.Range("K22").Value = 2.1
$ws.Range("T22").Value = 1.8
$ws.Range("W22").Value = 3.5
$ws.Range("X22").Value = 1.29
$ws.Range("AA22").Value = 1.8
$ws.Range("AB22").Value = 1.95
$ws.Range("AC22").Value = 9
$ws.Range("AG22").Value = 23
$ws.Range("AI22").Value = 10
$ws.Range("AM22").Value = 251
$ws.Range("AO22").Value = 11

# Row 23
$ws.Range("H23").Value = 3.7
$ws.Range("S23").Value = 2.05
$ws.Range("T23").Value = 1.8
$ws.Range("W23").Value = 3.5
$ws.Range("X23").Value = 1.29
$ws.Range("AD23").Value = 8.5
$ws.Range("AG23").Value = 15
$ws.Range("AI23").Value = 10
$ws.Range("AM23").Value = 351
$ws.Range("AP23").Value = 13

# Row 24
$ws.Range("S24").Value = 2.08
$ws.Range("T24").Value = 1.73

# Row 26
$ws.Range("K26").Value = 2.1
$ws.Range("M26").Value = 1.06
$ws.Range("N26").Value = 10
$ws.Range("O26").Value = 1.3
$ws.Range("P26").Value = 3.4
$ws.Range("S26").Value = 2
$ws.Range("T26").Value = 1.8
$ws.Range("W26").Value = 3.5
$ws.Range("X26").Value = 1.29
$ws.Range("AA26").Value = 1.8
$ws.Range("AB26").Value = 1.91
$ws.Range("AM26").Value = 251
$ws.Range("AO26").Value = 9.5

# Row 27
$ws.Range("G27").Value = 9
$ws.Range("H27").Value = 5.5
$ws.Range("I27").Value = 1.24
$ws.Range("J27").Value = 7.4
$ws.Range("K27").Value = 2.77
$ws.Range("L27").Value = 1.62
$ws.Range("M27").Value = 1.02
$ws.Range("N27").Value = 10
$ws.Range("O27").Value = 1.12
$ws.Range("P27").Value = 5.3
$ws.Range("S27").Value = 1.39
$ws.Range("T27").Value = 2.75
$ws.Range("W27").Value = 1.95
$ws.Range("X27").Value = 1.75
$ws.Range("Y27").Value = 1.23
$ws.Range("Z27").Value = 3.7
$ws.Range("AA27").Value = 1.78
$ws.Range("AB27").Value = 1.93
$ws.Range("AC27").Value = 32
$ws.Range("AE27").Value = 29
$ws.Range("AF27").Value = 250
$ws.Range("AG27").Value = 100
$ws.Range("AH27").Value = 75
$ws.Range("AI27").Value = 10
$ws.Range("AJ27").Value = 12
$ws.Range("AL27").Value = 75
$ws.Range("AM27").Value = 500
$ws.Range("AN27").Value = 10
$ws.Range("AO27").Value = 7.6
$ws.Range("AQ27").Value = 8.25
$ws.Range("AS27").Value = 23
